# Insert a new weekly price-report row for "Pepino ensalada" (Macroferia
# Regional de Talca) at row 379, pushing the existing rows 379..464 down to
# 380..465 (dimension grows from A1:R464 to A1:R465).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 379; this shifts rows 379-464
# down to 380-465, matching the rest of the sheet's existing data.
$ws.Rows.Item(379).Insert()

# Populate the newly inserted row 379 with the new record's data.
$ws.Cells.Item(379, 1).Value  = 5
$ws.Cells.Item(379, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(379, 3).Value  = "Maule"
$ws.Cells.Item(379, 4).Value  = 44785
$ws.Cells.Item(379, 5).Value  = 7
$ws.Cells.Item(379, 6).Value  = 100112043
$ws.Cells.Item(379, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(379, 8).Value  = "Sin especificar"
$ws.Cells.Item(379, 9).Value  = "Primera"
$ws.Cells.Item(379, 10).Value = 300
$ws.Cells.Item(379, 11).Value = 22000
$ws.Cells.Item(379, 12).Value = 22000
$ws.Cells.Item(379, 13).Value = 22000
$ws.Cells.Item(379, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(379, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(379, 16).Value = 367
$ws.Cells.Item(379, 17).Value = 60
$ws.Cells.Item(379, 18).Value = "Hortaliza"
